$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Revert "github for win" back to the numeric value 2222
$ws.Range("B5").Value = 2222

# Move the selection from B6 back to B5
$ws.Activate()
$ws.Range("B5").Select()
